$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Total TC Executed" value from 79 to 84
$ws.Range("D17").Value = 84

# Update "Total TC Passed" value from 5 to 0
$ws.Range("D21").Value = 0

# Update the selection/active cell on the sheet view
$ws.Range("D21:F24").Select()

$wb.Save()
